$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: demoQAHomepage / hyperlink to http://www.demoqa.com/
$ws.Range("A5").Value = "demoQAHomepage"
$ws.Range("B5").Value = "http://www.demoqa.com/"
$ws.Hyperlinks.Add($ws.Range("B5"), "http://www.demoqa.com/")

# Row 6: MoveAmountX / 50
$ws.Range("A6").Value = "MoveAmountX"
$ws.Range("B6").Value = 50
$ws.Range("B6").NumberFormat = "0"

# Row 7: MoveAmountY / 50
$ws.Range("A7").Value = "MoveAmountY"
$ws.Range("B7").Value = 50
$ws.Range("B7").NumberFormat = "0"

# Move the active selection to B7, matching the saved view state
$ws.Range("B7").Select()
